$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mood "None" option -------------------------------------------------
# Fill the Mood row (row 11, columns B:AF) with the new "None" value so it
# becomes a shared string and is used for every mood cell, mirroring the
# habit-tracker's existing placeholder pattern.
$ws.Range("B11:AF11").Value = "None"

# Extend the data validation list on that same range to include "None" as
# the first/default choice.
$dv = $ws.Range("B11:AF11").Validation
$dv.Modify(3, 1, 1, '"None, Happy, Neutral, Sad"')

# --- View / window state -------------------------------------------------
# Reset the scroll position back to the top-left corner (A1) instead of the
# far-right columns, reset the zoom back to 100%, and leave the selection on
# the newly edited cell B11.
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$win.Zoom = 100
$ws.Range("B11").Select()
